$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employees")

$ws.Cells.Item(2, 3).Value = 132
$ws.Cells.Item(2, 4).Value = 4000
$ws.Cells.Item(3, 3).Value = 37.75
$ws.Cells.Item(3, 4).Value = 1175
$ws.Cells.Item(4, 3).Value = 90
$ws.Cells.Item(4, 4).Value = 2710
$ws.Cells.Item(5, 3).Value = 75
$ws.Cells.Item(5, 4).Value = 2250
$ws.Cells.Item(6, 3).Value = 33.25
$ws.Cells.Item(6, 4).Value = 1020
$ws.Cells.Item(7, 3).Value = 73.75
$ws.Cells.Item(7, 4).Value = 2245
$ws.Cells.Item(8, 3).Value = 78.5
$ws.Cells.Item(8, 4).Value = 2375
$ws.Cells.Item(9, 3).Value = 76.5
$ws.Cells.Item(9, 4).Value = 2310
$ws.Cells.Item(10, 3).Value = 63.5
$ws.Cells.Item(10, 4).Value = 1905
$ws.Cells.Item(11, 3).Value = 55.5
$ws.Cells.Item(11, 4).Value = 1705
$ws.Cells.Item(12, 3).Value = 45.75
$ws.Cells.Item(12, 4).Value = 1390
$ws.Cells.Item(13, 3).Value = 104
$ws.Cells.Item(13, 4).Value = 3170
$ws.Cells.Item(14, 3).Value = 133.75
$ws.Cells.Item(14, 4).Value = 4045
$ws.Cells.Item(15, 3).Value = 69
$ws.Cells.Item(15, 4).Value = 2120
$ws.Cells.Item(16, 3).Value = 58
$ws.Cells.Item(16, 4).Value = 1755
$ws.Cells.Item(17, 3).Value = 66.5
$ws.Cells.Item(17, 4).Value = 2020
$ws.Cells.Item(18, 3).Value = 91
$ws.Cells.Item(18, 4).Value = 2730
$ws.Cells.Item(19, 3).Value = 64.75
$ws.Cells.Item(19, 4).Value = 1990
$ws.Cells.Item(20, 3).Value = 83.25
$ws.Cells.Item(20, 4).Value = 2540
$ws.Cells.Item(21, 3).Value = 76
$ws.Cells.Item(21, 4).Value = 2285
$ws.Cells.Item(22, 3).Value = 32.5
$ws.Cells.Item(22, 4).Value = 985
$ws.Cells.Item(23, 3).Value = 52.25
$ws.Cells.Item(23, 4).Value = 1590
$ws.Cells.Item(24, 3).Value = 68.5
$ws.Cells.Item(24, 4).Value = 2080
$ws.Cells.Item(25, 3).Value = 72.75
$ws.Cells.Item(25, 4).Value = 2225
$ws.Cells.Item(26, 3).Value = 54.5
$ws.Cells.Item(26, 4).Value = 1660
$ws.Cells.Item(27, 3).Value = 55.5
$ws.Cells.Item(27, 4).Value = 1675
$ws.Cells.Item(28, 3).Value = 68.75
$ws.Cells.Item(28, 4).Value = 2070
